$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($CellRef, $Value)
    $rng = $ws.Range($CellRef)
    $escaped = $Value.Replace('"', '""')
    $rng.Formula = '="' + $escaped + '"'
    $rng.Copy($rng) | Out-Null
    $rng.PasteSpecial(-4163) | Out-Null
}

Set-TextValue 'D2' '29.212.78'
Set-TextValue 'E2' '  -0.68%  '
Set-TextValue 'D3' '1.860.63'
Set-TextValue 'E3' '  -0.96%  '
Set-TextValue 'E4' '  -0.10%  '
Set-TextValue 'D5' '242.36'
Set-TextValue 'D6' '0.7024'
Set-TextValue 'E6' '  -2.06%  '
Set-TextValue 'D7' '0.9991'
Set-TextValue 'E7' '  -0.15%  '
Set-TextValue 'D8' '0.07839'
Set-TextValue 'E8' '  -1.65%  '
Set-TextValue 'D9' '0.3121'
Set-TextValue 'E9' '  -0.94%  '
Set-TextValue 'E10' '  -3.57%  '
Set-TextValue 'D11' '0.07805'
Set-TextValue 'E11' '  -3.94%  '
Set-TextValue 'D12' '1.865.93'
Set-TextValue 'E12' '  -0.93%  '
Set-TextValue 'D13' '5.148'
Set-TextValue 'E13' '  -1.78%  '
Set-TextValue 'D14' '92.90'
Set-TextValue 'E14' '  -2.55%  '
Set-TextValue 'D15' '0.6941'
Set-TextValue 'D16' '6.581'
Set-TextValue 'E16' '  +2.52%  '
Set-TextValue 'D17' '0.000008513'
Set-TextValue 'E17' '  +0.88%  '
Set-TextValue 'D18' '29.240.92'
Set-TextValue 'E18' '  -0.58%  '
Set-TextValue 'D19' '250.03'
Set-TextValue 'E19' '  -1.38%  '
Set-TextValue 'D20' '2.107.95'
Set-TextValue 'E20' '  -1.32%  '
Set-TextValue 'D21' '12.98'
Set-TextValue 'E21' '  -3.02%  '
Set-TextValue 'D22' '0.9987'
Set-TextValue 'E22' '  -0.14%  '
Set-TextValue 'D23' '7.636'
Set-TextValue 'E23' '  -0.59%  '
Set-TextValue 'E24' '  -0.07%  '
Set-TextValue 'D25' '0.1538'
Set-TextValue 'E25' '  -3.21%  '
Set-TextValue 'D26' '160.96'
Set-TextValue 'D27' '8.940'
Set-TextValue 'E27' '  -1.48%  '
Set-TextValue 'E28' '  -1.48%  '
Set-TextValue 'D29' '1.571'
Set-TextValue 'E29' '  +4.19%  '
Set-TextValue 'D30' '4.286'
Set-TextValue 'E30' '  -3.09%  '
Set-TextValue 'D31' '4.256'
Set-TextValue 'E31' '  -1.30%  '
Set-TextValue 'D32' '1.212'
Set-TextValue 'E32' '  -0.87%  '
Set-TextValue 'D33' '0.05248'
Set-TextValue 'E33' '  -1.51%  '
Set-TextValue 'D34' '0.7615'
Set-TextValue 'E34' '  +0.58%  '
Set-TextValue 'D35' '1.881'
Set-TextValue 'E35' '  -3.64%  '
Set-TextValue 'E36' '  +0.16%  '
Set-TextValue 'D37' '2.703'
Set-TextValue 'E37' '  +0.02%  '
Set-TextValue 'D38' '0.01864'
Set-TextValue 'E38' '  -1.69%  '
Set-TextValue 'D39' '1.231.74'
Set-TextValue 'E39' '  -3.32%  '
Set-TextValue 'D40' '2.724'
Set-TextValue 'E40' '  -1.48%  '
Set-TextValue 'D41' '0.9018'
Set-TextValue 'E41' '  -0.45%  '
Set-TextValue 'D42' '110.29'
Set-TextValue 'E42' '  -1.81%  '
Set-TextValue 'E43' '  -8.67%  '
Set-TextValue 'D44' '0.9981'
Set-TextValue 'E44' '  -0.26%  '
Set-TextValue 'D45' '68.18'
Set-TextValue 'E45' '  -8.31%  '
Set-TextValue 'D46' '2.004.87'
Set-TextValue 'E46' '  -1.32%  '
Set-TextValue 'E47' '  -4.12%  '
Set-TextValue 'D48' '0.5182'
Set-TextValue 'E48' '  -0.40%  '
Set-TextValue 'D49' '9.533'
Set-TextValue 'E49' '  +0.05%  '
Set-TextValue 'D50' '1.769'
Set-TextValue 'D51' '0.4266'

$excel.CutCopyMode = 0
